$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# ---------------------------------------------------------------------------
# 1) Underlying data edits that drive the recalculation cascade
# ---------------------------------------------------------------------------
$ws.Range("C9").Value = 38
$ws.Range("D9").Value = 6
$ws.Range("A9").Formula = "=AVERAGE(B9,568.3)"
$ws.Range("Q7").Value = 1.6243367818618635

# ---------------------------------------------------------------------------
# 2) New "S" column : 1/A^2
# ---------------------------------------------------------------------------
$ws.Range("S5").Formula = "=1/A5^2"
$ws.Range("S10").Formula = "=1/A10^2"
$ws.Range("S6:S9").Formula = "=1/A6^2"
$ws.Range("S11").Formula = "=1/A11^2"

# ---------------------------------------------------------------------------
# 3) New "J" column : minutes part of I
# ---------------------------------------------------------------------------
$ws.Range("J7").Formula = "=(I7-INT(I7))*60"
$ws.Range("J9:J11").Formula = "=(I9-INT(I9))*60"

# ---------------------------------------------------------------------------
# 4) Row 12 / Row 13 additions (NumberFormat applied before the J8 fill so
#    the new style entries land in the same order as the target workbook)
# ---------------------------------------------------------------------------
$ws.Range("I12").NumberFormat = "mmm-yy"

# J8 gets a plain white fill (no data)
$ws.Range("J8").Interior.ThemeColor = 2

$ws.Range("I13").Formula = "=11/60"
$ws.Range("J13").Formula = "=42/60"

# ---------------------------------------------------------------------------
# 5) Rows 30-34 : extra columns of angle data
# ---------------------------------------------------------------------------
$ws.Range("C30").Formula = "=289 + 26/60"
$ws.Range("D30").Value = 333

$ws.Range("C31").Formula = "=290 + 25/60"
$ws.Range("D31").Formula = "=332 + 4/60"
$ws.Range("E31").Formula = "=265.5 + 4/60"
$ws.Range("F31").Formula = "=356 + 10/60"

$ws.Range("C32").Formula = "=291 + 6/60"
$ws.Range("D32").Formula = "=331 + 21/60"

$ws.Range("C33").Formula = "=293 + 1/60"
$ws.Range("D33").Formula = "=329 + 24/60"

$ws.Range("C34").Formula = "=293.5 + 12/60"
$ws.Range("D34").Formula = "=328.5 + 11/60"

# ---------------------------------------------------------------------------
# 6) Column width for the new "S" column
# ---------------------------------------------------------------------------
$ws.Columns.Item(19).ColumnWidth = 12.71

# ---------------------------------------------------------------------------
# 7) Grab recalculated angle-of-minimum-deviation values for the new Sheet2
#    (these mirror column O for rows 5,7,9,10,11)
# ---------------------------------------------------------------------------
$o5 = $ws.Range("O5").Value()
$o7 = $ws.Range("O7").Value()
$o9 = $ws.Range("O9").Value()
$o10 = $ws.Range("O10").Value()
$o11 = $ws.Range("O11").Value()

# ---------------------------------------------------------------------------
# 8) New Sheet2
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws)
$ws2.Name = "Sheet2"
$ws2.Range("B2").Value = $o5
$ws2.Range("B3").Value = $o7
$ws2.Range("B4").Value = $o9
$ws2.Range("B5").Value = $o10
$ws2.Range("B6").Value = $o11
$ws2.Range("B1").Select()

# ---------------------------------------------------------------------------
# 9) Restore view state on Sheet1 (it must stay the active/visible tab)
# ---------------------------------------------------------------------------
$ws.Activate()
$excel.ActiveWindow.Zoom = 125
$excel.ActiveWindow.ScrollRow = 11
$ws.Range("D38").Select()

Write-Host "done"
